$d = $word.ActiveDocument
$t = $d.Tables(1)

# Each table row has a single cell whose sole paragraph/run holds the
# benchmark value. Replace the cell text wholesale (Word keeps the
# trailing cell-mark in place, and collapses any <w:tab/>-separated runs
# down to the single new run automatically).

function Set-CellText($rowIndex, $newText) {
    $cell = $t.Rows($rowIndex).Cells(1)
    $cell.Range.Text = $newText
}

Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "202"
Set-CellText 6 "0.00006"
Set-CellText 7 "0.00003"
Set-CellText 9 "0.00003"
Set-CellText 11 "0.00004"
Set-CellText 12 "0.00581"

Set-CellText 44 "100"
Set-CellText 45 "0.01"
Set-CellText 46 "154"
